$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, shifting existing rows 63-109 down to 64-110.
$ws.Rows("63").Insert()

# Populate the newly inserted row 63 with the new record's data.
$ws.Range("A63").Value = 5
$ws.Range("B63").Value = "Macroferia Regional de Talca"
$ws.Range("C63").Value = "Maule"
$ws.Range("D63").Value = 44603
$ws.Range("E63").Value = 7
$ws.Range("F63").Value = 100112030
$ws.Range("G63").Value = "Poroto granado"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 300
$ws.Range("K63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("M63").Value = 20000
$ws.Range("N63").Value = "`$/saco 25 kilos"
$ws.Range("O63").Value = "Región del Maule"
$ws.Range("P63").Value = 800
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
